$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("first_eval")

# Update existing rows 2-9, columns B:F with new values
$ws.Range("B2").Value = 0.6218823598496552
$ws.Range("C2").Value = 1.096605552064985
$ws.Range("D2").Value = 2.427837951941459
$ws.Range("E2").Value = 1.558152095253047
$ws.Range("F2").Value = 1.48260176597893

$ws.Range("B3").Value = 0.8174721314412605
$ws.Range("C3").Value = 1.10518048279258
$ws.Range("D3").Value = 2.163959630326707
$ws.Range("E3").Value = 1.471040322467983
$ws.Range("F3").Value = 1.27292596931789

$ws.Range("B4").Value = 0.9296915348393564
$ws.Range("C4").Value = 1.057526726501673
$ws.Range("D4").Value = 2.122174631855285
$ws.Range("E4").Value = 1.456768558095378
$ws.Range("F4").Value = 1.171408607494721

$ws.Range("B5").Value = 0.8447393623444465
$ws.Range("C5").Value = 0.9483023203837487
$ws.Range("D5").Value = 2.119205806994529
$ws.Range("E5").Value = 1.455749225311327
$ws.Range("F5").Value = 1.243456206856707

$ws.Range("B6").Value = 0.7998127184437969
$ws.Range("C6").Value = 0.8513425680011373
$ws.Range("D6").Value = 1.503415309836061
$ws.Range("E6").Value = 1.226138373037913
$ws.Range("F6").Value = 0.9796342431129901

$ws.Range("B7").Value = 0.9099748266321157
$ws.Range("C7").Value = 0.9718423566515457
$ws.Range("D7").Value = 1.455436698790223
$ws.Range("E7").Value = 1.206414812073452
$ws.Range("F7").Value = 0.8401222100961462
$ws.Range("G7").Value = 9

$ws.Range("B8").Value = 1.088978392985429
$ws.Range("C8").Value = 1.112266431153142
$ws.Range("D8").Value = 2.215538431729544
$ws.Range("E8").Value = 1.488468485299418
$ws.Range("F8").Value = 1.111574284341132
$ws.Range("G8").Value = 6

$ws.Range("B9").Value = 1.309335564656196
$ws.Range("C9").Value = 1.309335564656196
$ws.Range("D9").Value = 2.208706261879395
$ws.Range("E9").Value = 1.486171679813404
$ws.Range("F9").Value = 0.8611155331944436
$ws.Range("G9").Value = 3

# New row 10 for Q8
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A10").Value = "Q8"
$ws.Range("B10").Value = 0.5344944032044353
$ws.Range("C10").Value = 0.5344944032044353
$ws.Range("D10").Value = 0.2856842670568654
$ws.Range("E10").Value = 0.5344944032044353
$ws.Range("G10").Value = 1
